$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 already contains a record identical to the new one that needs
# to be appended (same timestamp/value/N-A pattern), so copy it down to
# row 12 to add the new data row reported by Adafruit IO.
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4163)
$excel.CutCopyMode = 0
